# Add a new "2022" column (S) to the SDG indicator table, mirroring the
# existing "2021" column (R) formatting, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column R (rows 2-6) into column S so the new
# cells inherit the same styles (thick-bottom border row, header style,
# integer style, one-decimal styles, etc.) as the rest of the table.
$ws.Range("R2:R6").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New "2022" column values.
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Move the selection, as recorded in the saved workbook view.
$ws.Range("C19").Select() | Out-Null
